# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Membrillo" (Vega Central Mapocho de
# Santiago) at row 7, pushing the existing rows (old 7..37) down to 9..39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before the current row 7 (old row 7 -> row 9, ... old row 37 -> row 39)
$ws.Range("A7:A8").EntireRow.Insert()

# Columns that are constant across every record in this sheet
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100104
$producto  = "Frutos de pepita"
$categoriaId = 100104003
$categoria = "Membrillo"
$variedad  = "Champion"

function Set-MembrilloRow($row, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $origen, $pkg, $kg) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $pmin
    $ws.Cells.Item($row, 15).Value = $pmax
    $ws.Cells.Item($row, 16).Value = $pprom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $pkg
    $ws.Cells.Item($row, 20).Value = $kg
}

# New row 7: 2022-04-19, Primera, 220 @ $18000 ($/caja 18 kilos granel), Región de O'Higgins
Set-MembrilloRow 7 44670 "Primera" 220 18000 18000 18000 "$/caja 18 kilos granel" "Región de O'Higgins" 1000 18

# New row 8: 2022-04-19, Segunda, 250 @ $14400 ($/caja 18 kilos granel), Región de O'Higgins
Set-MembrilloRow 8 44670 "Segunda" 250 14400 14400 14400 "$/caja 18 kilos granel" "Región de O'Higgins" 800 18

# Make sure the date cells keep the date number format used by the rest of column D
$ws.Range("D7:D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
